$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 553.17645
$ws.Range("I55").Value = 231.66667
$ws.Range("J55").Value = 728.5454999999999
$ws.Range("K55").Value = 231.66667
$ws.Range("L55").Value = 728.5454999999999
$ws.Range("M55").Value = -17.66667000000001
$ws.Range("N55").Value = -1156.5455

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2306.625
$ws.Range("I88").Value = 656.3333
$ws.Range("J88").Value = 3296.8
$ws.Range("K88").Value = 656.3333
$ws.Range("L88").Value = 3296.8
$ws.Range("M88").Value = -250.3333
$ws.Range("N88").Value = -4108.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2306.625
$ws.Range("I91").Value = 656.3333
$ws.Range("J91").Value = 3296.8
$ws.Range("K91").Value = 656.3333
$ws.Range("L91").Value = 3296.8
$ws.Range("M91").Value = 747.6667
$ws.Range("N91").Value = -6104.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3254.68
$ws.Range("J137").Value = 4284.643
$ws.Range("L137").Value = 12853.929
$ws.Range("N137").Value = -17953.929

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2859.6667
$ws.Range("I141").Value = 431.6
$ws.Range("K141").Value = 1294.8
$ws.Range("M141").Value = 3885.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1527.4
$ws.Range("I61").Value = 1540.0714
$ws.Range("K61").Value = 1540.0714
$ws.Range("M61").Value = -1328.0714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1878.5714
$ws.Range("I102").Value = 1694.125
$ws.Range("K102").Value = 1694.125
$ws.Range("M102").Value = -72.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6003
$ws.Range("I110").Value = 6003
$ws.Range("K110").Value = 6003
$ws.Range("M110").Value = -3958

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2917.739
$ws.Range("I122").Value = 2133.1875
$ws.Range("K122").Value = 6399.5625
$ws.Range("M122").Value = -3949.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1527.4
$ws.Range("I136").Value = 1540.0714
$ws.Range("K136").Value = 4620.2142
$ws.Range("M136").Value = -2070.2142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 20000
$ws.Range("K38").Value = 20000
$ws.Range("M38").Value = -19584

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 580.46155
$ws.Range("I94").Value = 553.8333
$ws.Range("K94").Value = 553.8333
$ws.Range("M94").Value = -102.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1722.9259
$ws.Range("I107").Value = 1322.25
$ws.Range("K107").Value = 1322.25
$ws.Range("M107").Value = 597.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2222
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2222
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4808.7
$ws.Range("I31").Value = 3272.8572
$ws.Range("K31").Value = 3272.8572
$ws.Range("M31").Value = -2977.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4808.7
$ws.Range("I34").Value = 3272.8572
$ws.Range("K34").Value = 3272.8572
$ws.Range("M34").Value = -3070.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 15001.5
$ws.Range("J86").Value = 16144.571
$ws.Range("L86").Value = 16144.571
$ws.Range("N86").Value = -18390.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 15001.5
$ws.Range("J89").Value = 16144.571
$ws.Range("L89").Value = 80722.855
$ws.Range("N89").Value = -91954.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3559.1538
$ws.Range("I105").Value = 1819.7142
$ws.Range("K105").Value = 1819.7142
$ws.Range("M105").Value = -72.71419999999989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1144.3334
$ws.Range("I107").Value = 821.4
$ws.Range("K107").Value = 821.4
$ws.Range("M107").Value = 1098.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2222
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2222
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4292.2856
$ws.Range("I134").Value = 3592.4167
$ws.Range("K134").Value = 10777.2501
$ws.Range("M134").Value = -8242.250100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 7474.5
$ws.Range("I59").Value = 7474.5
$ws.Range("K59").Value = 22423.5
$ws.Range("M59").Value = -21883.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2185.5
$ws.Range("I81").Value = 600
$ws.Range("K81").Value = 1800
$ws.Range("M81").Value = -677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2185.5
$ws.Range("I84").Value = 600
$ws.Range("K84").Value = 5400
$ws.Range("M84").Value = 216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 550
$ws.Range("I103").Value = 325
$ws.Range("K103").Value = 975
$ws.Range("M103").Value = -96

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1350.5
$ws.Range("I107").Value = 1025.2
$ws.Range("J107").Value = 1582.8572
$ws.Range("K107").Value = 3075.6
$ws.Range("L107").Value = 4748.571599999999
$ws.Range("M107").Value = -1155.6
$ws.Range("N107").Value = -8588.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2055.4
$ws.Range("I113").Value = 4899
$ws.Range("K113").Value = 14697
$ws.Range("M113").Value = -12527

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 12499
$ws.Range("J139").Value = 12499
$ws.Range("L139").Value = 37497
$ws.Range("N139").Value = -47777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2901.2144
$ws.Range("I140").Value = 2355.1538
$ws.Range("J140").Value = 10000
$ws.Range("K140").Value = 7065.4614
$ws.Range("L140").Value = 30000
$ws.Range("M140").Value = -1885.4614
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 6714.2856
$ws.Range("I92").Value = 25000
$ws.Range("J92").Value = 3666.6667
$ws.Range("K92").Value = 25000
$ws.Range("L92").Value = 3666.6667
$ws.Range("M92").Value = -23128
$ws.Range("N92").Value = -7410.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1754.75
$ws.Range("I102").Value = 769
$ws.Range("K102").Value = 769
$ws.Range("M102").Value = 853

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 410606.66
$ws.Range("I122").Value = 93109
$ws.Range("J122").Value = 628886.3
$ws.Range("K122").Value = 279327
$ws.Range("L122").Value = 1886658.9
$ws.Range("M122").Value = -276877
$ws.Range("N122").Value = -1891558.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6695.476
$ws.Range("I62").Value = 3666.3333
$ws.Range("J62").Value = 7200.3335
$ws.Range("K62").Value = 3666.3333
$ws.Range("L62").Value = 7200.3335
$ws.Range("M62").Value = -3042.3333
$ws.Range("N62").Value = -8448.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6695.476
$ws.Range("I65").Value = 3666.3333
$ws.Range("J65").Value = 7200.3335
$ws.Range("K65").Value = 18331.6665
$ws.Range("L65").Value = 36001.6675
$ws.Range("M65").Value = -15211.6665
$ws.Range("N65").Value = -42241.6675
